$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update relabeled subcategory text values in column H
$ws.Range("H6").Value = "bar chart(s)"
$ws.Range("H9").Value = "line graph(s)"
$ws.Range("H14").Value = "data collection, data analysis, data gathering diagram"
$ws.Range("H15").Value = "bar chart(s)"
$ws.Range("H20").Value = "line graph(s)"

# Remove the "is_viewed" column entirely (column I)
$ws.Range("I1:I26").EntireColumn.Delete()
